$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the file name reference from the PNG screenshot to the PDF export
$ws.Range("B2").Value = "Data Model Links.pdf"

# Move the active selection to C9 (matches the saved cursor position)
$ws.Range("C9").Select()
